# Apply the "Add files via upload" edit to the Orlando attractions sheet:
#  - Clear the one-off large/colored font that had been applied to A2
#    ("Despicable Me Minion Mayhem"), restoring it to the sheet's default
#    style (and letting the row reclaim its default height).
#  - Append three new attraction rows (name, lat, lon) after the existing
#    eight rows of data.
#  - Leave the active selection on F8, matching the saved workbook view.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Revert the bespoke formatting on A2 back to the workbook default ---
$ws.Range("A2").Style = "Normal"
$ws.Rows.Item(2).AutoFit()

# --- Append the new attractions ---
$ws.Range("A9").Value = "Revenge of the Mummy"
$ws.Range("B9").Value = 28.4769
$ws.Range("C9").Value = -81.469886

$ws.Range("A10").Value = "Transformers: The Ride 3D"
$ws.Range("B10").Value = 28.476361
$ws.Range("C10").Value = -81.468364

$ws.Range("A11").Value = "E.T. Adventure"
$ws.Range("B11").Value = 28.4776
$ws.Range("C11").Value = -81.4665

# --- Match the saved selection state ---
$ws.Range("F8").Select() | Out-Null
